$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 64, pushing every
# following row (old 64..211) down by one (to 65..212). Excel's row
# insert handles the shifting (and the dimension/used-range update)
# automatically; only the brand-new row needs its values populated.
$ws.Rows("64:64").Insert()

$ws.Cells.Item(64, 1).Value = 4
$ws.Cells.Item(64, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(64, 3).Value = "Los Lagos"
$ws.Cells.Item(64, 4).Value = 44622
$ws.Cells.Item(64, 5).Value = 10
$ws.Cells.Item(64, 6).Value = 100112017
$ws.Cells.Item(64, 7).Value = "Apio"
$ws.Cells.Item(64, 8).Value = "Americana (o)"
$ws.Cells.Item(64, 9).Value = "Primera"
$ws.Cells.Item(64, 10).Value = 25
$ws.Cells.Item(64, 11).Value = 13000
$ws.Cells.Item(64, 12).Value = 13000
$ws.Cells.Item(64, 13).Value = 13000
$ws.Cells.Item(64, 14).Value = '$/docena de matas'
$ws.Cells.Item(64, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(64, 16).Value = 2167
$ws.Cells.Item(64, 17).Value = 6
$ws.Cells.Item(64, 18).Value = "Hortaliza"
